$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Containers & Packaging(12)', 0.7001438813779158),
    @('Mortgage Real Estate Investment Trust...(16)', 0.6904882631932447),
    @('Equity Real Estate Investment Trusts ...(98)', 0.6398474642353622),
    @('Road & Rail(22)', 0.6293722548865807),
    @('ETF(303)', 0.6139360250295423),
    @('Multi-Utilities(18)', 0.5813759527450517),
    @('Energy Equipment & Services(38)', 0.561435902338057),
    @('Building Products(24)', 0.5595365801594221),
    @('Airlines(14)', 0.5044110339867112),
    @('Life Sciences Tools & Services(19)', 0.4932196026268766),
    @('Machinery(86)', 0.4899199090905587),
    @('Internet & Direct Marketing Retail(15)', 0.4756808738508365),
    @('Construction & Engineering(21)', 0.4678149584645998),
    @('Trading Companies & Distributors(25)', 0.4577437672919444),
    @('Gas Utilities(12)', 0.4484509301746201),
    @('Semiconductors & Semiconductor Equipment(70)', 0.4414651210783485),
    @('Auto Components(21)', 0.4372359952101255),
    @('Consumer Finance(15)', 0.4372342122136115),
    @('Banks(251)', 0.4284718710030003),
    @('Leisure Products(11)', 0.4237333247381553),
    @('Electrical Equipment(28)', 0.4227163602566927),
    @('Household Products(10)', 0.4055378415224902),
    @('Real Estate Management & Development(23)', 0.3834147971440759),
    @('Chemicals(52)', 0.3687769932485391),
    @('Air Freight & Logistics(11)', 0.3677915872076659),
    @('Oil, Gas & Consumable Fuels(125)', 0.3563361837970567),
    @('Technology Hardware, Storage & Periph...(19)', 0.3543194975450485),
    @('IT Services(52)', 0.3391849322854577),
    @('Electronic Equipment, Instruments & C...(78)', 0.3379169324604882),
    @('Commercial Services & Supplies(52)', 0.3308361528058588),
    @('Software(70)', 0.3287416546126967),
    @('Health Care Providers & Services(47)', 0.326154717829893),
    @('Electric Utilities(28)', 0.324950483858519),
    @('Hotels, Restaurants & Leisure(51)', 0.324600144318076),
    @('Capital Markets(76)', 0.3204853498788125),
    @('Household Durables(39)', 0.3181525400969379),
    @('Beverages(21)', 0.3127055360960744),
    @('Aerospace & Defense(37)', 0.3087510583339029),
    @('Insurance(75)', 0.3024132463637215),
    @('Professional Services(35)', 0.2985858803840739),
    @('Thrifts & Mortgage Finance(47)', 0.2939761905456659),
    @('Marine(15)', 0.2834946746149619),
    @('Textiles, Apparel & Luxury Goods(29)', 0.2794280827201666),
    @('Health Care Equipment & Supplies(86)', 0.2756420031878652),
    @('Metals & Mining(106)', 0.2571860721523278),
    @('Food Products(46)', 0.2520102537249483),
    @('Biotechnology(128)', 0.2425934687722011),
    @('Media(42)', 0.2372918794920321),
    @('Communications Equipment(45)', 0.2308892914139308),
    @('Specialty Retail(59)', 0.2129698478940226),
    @('Entertainment(22)', 0.2067035958370407),
    @('Pharmaceuticals(53)', 0.1629597184251727)
)

$rowCount = $data.Count
$startRow = 2
$endRow = $startRow + $rowCount - 1

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
